# CryCompanywiseStockReport_1.xlsx - restore swapped row pairs.
#
# A set of adjacent row pairs in the report had their data (everything
# except the running serial number in column A) swapped between the two
# rows. This script swaps columns B:G (Item Code, Item Name, Rate, MRP,
# Qty, Value) back between each pair of rows, leaving column A (and the
# empty H:M columns) untouched.
#
# NOTE: reading/writing the COM `.Value` property on this host returns the
# property descriptor string instead of the live cell value, so `.Value2`
# is used instead (behaves like `.Value` for numbers/strings, without the
# date-as-Variant(VT_DATE) wrapping we don't need here).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-RowRange($row1, $row2, $firstCol, $lastCol) {
    for ($col = $firstCol; $col -le $lastCol; $col++) {
        $cell1 = $ws.Cells.Item($row1, $col)
        $cell2 = $ws.Cells.Item($row2, $col)
        $value1 = $cell1.Value2
        $value2 = $cell2.Value2
        $cell1.Value2 = $value2
        $cell2.Value2 = $value1
    }
}

# (row1, row2) pairs whose B:G (cols 2..7) content needs to be swapped back.
$rowPairs = @(
    @(183, 184),
    @(264, 265),
    @(279, 280),
    @(313, 314),
    @(317, 318),
    @(350, 351),
    @(355, 356),
    @(375, 376),
    @(379, 380),
    @(382, 383),
    @(421, 422),
    @(431, 432),
    @(536, 537),
    @(579, 580),
    @(590, 591),
    @(599, 600),
    @(687, 688),
    @(720, 721),
    @(872, 873)
)

foreach ($pair in $rowPairs) {
    Swap-RowRange $pair[0] $pair[1] 2 7
}
